$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells stay text (avoid Excel auto-converting numeric-looking
# strings like "1.00" or "595.28" into real numbers and losing formatting/precision).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.484.39'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.750.81'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.28'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.33'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.746.73'
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.49'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000259'
$ws.Range('E13').Value = '  -5.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.16'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.380.46'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.738.17'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.465.36'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.89'
$ws.Range('E18').Value = '  -3.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.02'
$ws.Range('E19').Value = '  -2.08%  '
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '467.75'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.19'
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.20'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.03'
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.14'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.895.58'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.78'
$ws.Range('E31').Value = '  -4.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.30'
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.89'
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('E34').Value = '  -1.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.25'
$ws.Range('E35').Value = '  +1.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.704.97'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  -2.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.40'
$ws.Range('E39').Value = '  -10.55%  '
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.80'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.305'
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.59'
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.14'
$ws.Range('E47').Value = '  +11.21%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.93'
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.81'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '146.21'
$ws.Range('E50').Value = '  +5.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '391.82'
$ws.Range('E51').Value = '  -0.82%  '
